$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Range("B1").Value = "groomed_file"
$ws.Range("B2").Value = "sphere10_DT.nrrd"

$ws.Range("C1").Value = "transform_file"
$ws.Range("C2").Value = " 1.000000 0.000000 0.000000 0.000000 1.000000 0.000000 0.000000 0.000000 1.000000 -0.002400 -0.002400 -0.002400"

$ws.Range("B3").Value = "sphere20_DT.nrrd"
$ws.Range("C3").Value = " 1.000000 0.000000 0.000000 0.000000 1.000000 0.000000 0.000000 0.000000 1.000000 -0.000599 -0.000599 -0.000599"

$ws.Range("B4").Value = "sphere30_DT.nrrd"
$ws.Range("C4").Value = " 1.000000 0.000000 0.000000 0.000000 1.000000 0.000000 0.000000 0.000000 1.000000 -0.000265 -0.000265 -0.000265"

$ws.Range("B5").Value = "sphere40_DT.nrrd"
$ws.Range("C5").Value = " 1.000000 0.000000 0.000000 0.000000 1.000000 0.000000 0.000000 0.000000 1.000000 -0.000149 -0.000149 -0.000149"
